$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.8733400000000001
$ws.Range("H2").Value = 2.62002
$ws.Range("I2").Value = 0.2319025556440181
$ws.Range("J2").Value = 0.2319025556440181
$ws.Range("M2").Value = 38.45264233333334
$ws.Range("N2").Value = 115.357927
$ws.Range("O2").Value = 0.2975040117664333
$ws.Range("P2").Value = 0.2975040117664332
$ws.Range("Q2").Value = 33.58223065539334
$ws.Range("R2").Value = 302.2400758985401
$ws.Range("S2").Value = 0.06899194064298392
$ws.Range("T2").Value = 0.0689919406429839

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.8733400000000001
$ws.Range("H3").Value = 2.62002
$ws.Range("I3").Value = 0.2319025556440181
$ws.Range("J3").Value = 0.2319025556440181
$ws.Range("O3").Value = 0.3694391181876273
$ws.Range("P3").Value = 0.3694391181876272
$ws.Range("Q3").Value = 41.70226010210001
$ws.Range("R3").Value = 375.3203409189001
$ws.Range("S3").Value = 0.08567387566258322
$ws.Range("T3").Value = 0.08567387566258321

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.8733400000000001
$ws.Range("H4").Value = 2.62002
$ws.Range("I4").Value = 0.2319025556440181
$ws.Range("J4").Value = 0.2319025556440181
$ws.Range("M4").Value = 18.63107466666667
$ws.Range("N4").Value = 55.893224
$ws.Range("O4").Value = 0.1441466469015163
$ws.Range("P4").Value = 0.1441466469015162
$ws.Range("Q4").Value = 16.27126274938667
$ws.Range("R4").Value = 146.44136474448
$ws.Range("S4").Value = 0.03342797580397751
$ws.Range("T4").Value = 0.0334279758039775

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.8733400000000001
$ws.Range("H5").Value = 2.62002
$ws.Range("I5").Value = 0.2319025556440181
$ws.Range("J5").Value = 0.2319025556440181
$ws.Range("M5").Value = 24.41680433333333
$ws.Range("N5").Value = 73.25041299999999
$ws.Range("O5").Value = 0.1889102231444233
$ws.Range("P5").Value = 0.1889102231444233
$ws.Range("Q5").Value = 21.32417189647333
$ws.Range("R5").Value = 191.91754706826
$ws.Range("S5").Value = 0.0438087635344735
$ws.Range("T5").Value = 0.0438087635344735

# Row 6
$ws.Range("I6").Value = 0.07113291057171066
$ws.Range("J6").Value = 0.07113291057171067
$ws.Range("M6").Value = 38.45264233333334
$ws.Range("N6").Value = 115.357927
$ws.Range("O6").Value = 0.2975040117664333
$ws.Range("P6").Value = 0.2975040117664332
$ws.Range("Q6").Value = 10.300886091465
$ws.Range("R6").Value = 92.707974823185
$ws.Range("S6").Value = 0.02116232626370685
$ws.Range("T6").Value = 0.02116232626370685

# Row 7
$ws.Range("I7").Value = 0.07113291057171066
$ws.Range("J7").Value = 0.07113291057171067
$ws.Range("O7").Value = 0.3694391181876273
$ws.Range("P7").Value = 0.3694391181876272
$ws.Range("S7").Value = 0.02627927975573214
$ws.Range("T7").Value = 0.02627927975573214

# Row 8
$ws.Range("I8").Value = 0.07113291057171066
$ws.Range("J8").Value = 0.07113291057171067
$ws.Range("M8").Value = 18.63107466666667
$ws.Range("N8").Value = 55.893224
$ws.Range("O8").Value = 0.1441466469015163
$ws.Range("P8").Value = 0.1441466469015162
$ws.Range("Q8").Value = 4.99098543708
$ws.Range("R8").Value = 44.91886893372001
$ws.Range("S8").Value = 0.01025357054325751
$ws.Range("T8").Value = 0.01025357054325751

# Row 9
$ws.Range("I9").Value = 0.07113291057171066
$ws.Range("J9").Value = 0.07113291057171067
$ws.Range("M9").Value = 24.41680433333333
$ws.Range("N9").Value = 73.25041299999999
$ws.Range("O9").Value = 0.1889102231444233
$ws.Range("P9").Value = 0.1889102231444233
$ws.Range("Q9").Value = 6.540895628834999
$ws.Range("R9").Value = 58.868060659515
$ws.Range("S9").Value = 0.01343773400901417
$ws.Range("T9").Value = 0.01343773400901417

# Row 10
$ws.Range("G10").Value = 0.06721833333333334
$ws.Range("H10").Value = 0.201655
$ws.Range("I10").Value = 0.01784883697773089
$ws.Range("J10").Value = 0.01784883697773089
$ws.Range("M10").Value = 38.45264233333334
$ws.Range("N10").Value = 115.357927
$ws.Range("O10").Value = 0.2975040117664333
$ws.Range("P10").Value = 0.2975040117664332
$ws.Range("Q10").Value = 2.584722529909445
$ws.Range("R10").Value = 23.262502769185
$ws.Range("S10").Value = 0.005310100606239998
$ws.Range("T10").Value = 0.005310100606239997

# Row 11
$ws.Range("G11").Value = 0.06721833333333334
$ws.Range("H11").Value = 0.201655
$ws.Range("I11").Value = 0.01784883697773089
$ws.Range("J11").Value = 0.01784883697773089
$ws.Range("O11").Value = 0.3694391181876273
$ws.Range("P11").Value = 0.3694391181876272
$ws.Range("Q11").Value = 3.209696590441667
$ws.Range("R11").Value = 28.887269313975
$ws.Range("S11").Value = 0.006594058593727612
$ws.Range("T11").Value = 0.006594058593727611

# Row 12
$ws.Range("G12").Value = 0.06721833333333334
$ws.Range("H12").Value = 0.201655
$ws.Range("I12").Value = 0.01784883697773089
$ws.Range("J12").Value = 0.01784883697773089
$ws.Range("M12").Value = 18.63107466666667
$ws.Range("N12").Value = 55.893224
$ws.Range("O12").Value = 0.1441466469015163
$ws.Range("P12").Value = 0.1441466469015162
$ws.Range("Q12").Value = 1.252349787302222
$ws.Range("R12").Value = 11.27114808572
$ws.Range("S12").Value = 0.002572850001431701
$ws.Range("T12").Value = 0.0025728500014317

# Row 13
$ws.Range("G13").Value = 0.06721833333333334
$ws.Range("H13").Value = 0.201655
$ws.Range("I13").Value = 0.01784883697773089
$ws.Range("J13").Value = 0.01784883697773089
$ws.Range("M13").Value = 24.41680433333333
$ws.Range("N13").Value = 73.25041299999999
$ws.Range("O13").Value = 0.1889102231444233
$ws.Range("P13").Value = 0.1889102231444233
$ws.Range("Q13").Value = 1.641256892612778
$ws.Range("R13").Value = 14.771312033515
$ws.Range("S13").Value = 0.003371827776331575
$ws.Range("T13").Value = 0.003371827776331575

# Row 14
$ws.Range("G14").Value = 2.557535
$ws.Range("H14").Value = 7.672605
$ws.Range("I14").Value = 0.6791156968065403
$ws.Range("J14").Value = 0.6791156968065403
$ws.Range("M14").Value = 38.45264233333334
$ws.Range("N14").Value = 115.357927
$ws.Range("O14").Value = 0.2975040117664333
$ws.Range("P14").Value = 0.2975040117664332
$ws.Range("Q14").Value = 98.34397860998168
$ws.Range("R14").Value = 885.095807489835
$ws.Range("S14").Value = 0.2020396442535025
$ws.Range("T14").Value = 0.2020396442535024

# Row 15
$ws.Range("G15").Value = 2.557535
$ws.Range("H15").Value = 7.672605
$ws.Range("I15").Value = 0.6791156968065403
$ws.Range("J15").Value = 0.6791156968065403
$ws.Range("O15").Value = 0.3694391181876273
$ws.Range("P15").Value = 0.3694391181876272
$ws.Range("Q15").Value = 122.123101873525
$ws.Range("R15").Value = 1099.107916861725
$ws.Range("S15").Value = 0.2508919041755843
$ws.Range("T15").Value = 0.2508919041755842

# Row 16
$ws.Range("G16").Value = 2.557535
$ws.Range("H16").Value = 7.672605
$ws.Range("I16").Value = 0.6791156968065403
$ws.Range("J16").Value = 0.6791156968065403
$ws.Range("M16").Value = 18.63107466666667
$ws.Range("N16").Value = 55.893224
$ws.Range("O16").Value = 0.1441466469015163
$ws.Range("P16").Value = 0.1441466469015162
$ws.Range("Q16").Value = 47.64962554761333
$ws.Range("R16").Value = 428.84662992852
$ws.Range("S16").Value = 0.09789225055284954
$ws.Range("T16").Value = 0.09789225055284952

# Row 17
$ws.Range("G17").Value = 2.557535
$ws.Range("H17").Value = 7.672605
$ws.Range("I17").Value = 0.6791156968065403
$ws.Range("J17").Value = 0.6791156968065403
$ws.Range("M17").Value = 24.41680433333333
$ws.Range("N17").Value = 73.25041299999999
$ws.Range("O17").Value = 0.1889102231444233
$ws.Range("P17").Value = 0.1889102231444233
$ws.Range("Q17").Value = 62.44683167065166
$ws.Range("R17").Value = 562.021485035865
$ws.Range("S17").Value = 0.128291897824604
$ws.Range("T17").Value = 0.128291897824604
